# Update Leve profit-tracking values (scheduled market-price refresh)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (46 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2361.7307
$ws.Range("I40").Value = 2288.5293
$ws.Range("K40").Value = 2288.5293
$ws.Range("M40").Value = -2113.5293
$ws.Range("H100").Value = 37753.18
$ws.Range("I100").Value = 43639.168
$ws.Range("K100").Value = 43639.168
$ws.Range("M100").Value = -43098.168
$ws.Range("H106").Value = 7197.75
$ws.Range("I106").Value = 7197.75
$ws.Range("K106").Value = 7197.75
$ws.Range("M106").Value = -6566.75
$ws.Range("H113").Value = 3712
$ws.Range("I113").Value = 2899.5
$ws.Range("K113").Value = 2899.5
$ws.Range("M113").Value = 354.5
$ws.Range("H116").Value = 11938.3
$ws.Range("I116").Value = 11749.333
$ws.Range("K116").Value = 11749.333
$ws.Range("M116").Value = -8307.333000000001
$ws.Range("H118").Value = 729.0833
$ws.Range("I118").Value = 752.1818
$ws.Range("K118").Value = 2256.5454
$ws.Range("M118").Value = -599.5454
$ws.Range("H132").Value = 1816.4706
$ws.Range("I132").Value = 1730.6923
$ws.Range("J132").Value = 2095.25
$ws.Range("K132").Value = 5192.0769
$ws.Range("L132").Value = 6285.75
$ws.Range("M132").Value = -2662.0769
$ws.Range("N132").Value = -11345.75
$ws.Range("H135").Value = 1556.0834
$ws.Range("I135").Value = 1668
$ws.Range("K135").Value = 15012
$ws.Range("M135").Value = -12477
$ws.Range("H137").Value = 7023.727
$ws.Range("I137").Value = 2954.8667
$ws.Range("J137").Value = 10414.444
$ws.Range("K137").Value = 8864.6001
$ws.Range("L137").Value = 31243.332
$ws.Range("M137").Value = -6314.6001
$ws.Range("N137").Value = -36343.33199999999
$ws.Range("H141").Value = 3355.7778
$ws.Range("J141").Value = 3068
$ws.Range("L141").Value = 9204
$ws.Range("N141").Value = -19564

# --- Sheet: ARM (34 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 16367.375
$ws.Range("I45").Value = 16367.375
$ws.Range("K45").Value = 16367.375
$ws.Range("M45").Value = -15990.375
$ws.Range("H61").Value = 3352.0833
$ws.Range("J61").Value = 2633
$ws.Range("L61").Value = 2633
$ws.Range("N61").Value = -3057
$ws.Range("H88").Value = 899
$ws.Range("I88").Value = 935.875
$ws.Range("J88").Value = 825.25
$ws.Range("K88").Value = 935.875
$ws.Range("L88").Value = 825.25
$ws.Range("M88").Value = -529.875
$ws.Range("N88").Value = -1637.25
$ws.Range("H91").Value = 899
$ws.Range("I91").Value = 935.875
$ws.Range("J91").Value = 825.25
$ws.Range("K91").Value = 935.875
$ws.Range("L91").Value = 825.25
$ws.Range("M91").Value = 468.125
$ws.Range("N91").Value = -3633.25
$ws.Range("H132").Value = 6650.1763
$ws.Range("I132").Value = 6650.1763
$ws.Range("K132").Value = 19950.5289
$ws.Range("M132").Value = -17420.5289
$ws.Range("H136").Value = 3352.0833
$ws.Range("J136").Value = 2633
$ws.Range("L136").Value = 7899
$ws.Range("N136").Value = -12999
$ws.Range("H139").Value = 99698
$ws.Range("J139").Value = 99698
$ws.Range("L139").Value = 99698
$ws.Range("N139").Value = -109978

# --- Sheet: CRP (40 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2387.65
$ws.Range("I31").Value = 1716.4375
$ws.Range("K31").Value = 1716.4375
$ws.Range("M31").Value = -1421.4375
$ws.Range("H34").Value = 2387.65
$ws.Range("I34").Value = 1716.4375
$ws.Range("K34").Value = 1716.4375
$ws.Range("M34").Value = -1514.4375
$ws.Range("H99").Value = 5786.5
$ws.Range("I99").Value = 5490.6665
$ws.Range("J99").Value = 6082.3335
$ws.Range("K99").Value = 5490.6665
$ws.Range("L99").Value = 6082.3335
$ws.Range("M99").Value = -3992.6665
$ws.Range("N99").Value = -9078.333500000001
$ws.Range("H122").Value = 1631.5454
$ws.Range("I122").Value = 1438.5555
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4315.666499999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1865.666499999999
$ws.Range("N122").Value = -12400
$ws.Range("H126").Value = 5786.5
$ws.Range("I126").Value = 5490.6665
$ws.Range("J126").Value = 6082.3335
$ws.Range("K126").Value = 16471.9995
$ws.Range("L126").Value = 18247.0005
$ws.Range("M126").Value = -14001.9995
$ws.Range("N126").Value = -23187.0005
$ws.Range("H132").Value = 7801.9697
$ws.Range("I132").Value = 4356.3447
$ws.Range("K132").Value = 13069.0341
$ws.Range("M132").Value = -10539.0341
$ws.Range("H134").Value = 5411.8887
$ws.Range("I134").Value = 5934.3335
$ws.Range("J134").Value = 2799.6667
$ws.Range("K134").Value = 17803.0005
$ws.Range("L134").Value = 8399.000100000001
$ws.Range("M134").Value = -15268.0005
$ws.Range("N134").Value = -13469.0001

# --- Sheet: CUL (48 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1985.0968
$ws.Range("I5").Value = 1103.7
$ws.Range("J5").Value = 2404.8096
$ws.Range("K5").Value = 3311.1
$ws.Range("L5").Value = 7214.4288
$ws.Range("M5").Value = -3199.1
$ws.Range("N5").Value = -7438.4288
$ws.Range("H12").Value = 570.4737
$ws.Range("I12").Value = 52.714287
$ws.Range("K12").Value = 158.142861
$ws.Range("M12").Value = 14.85713900000002
$ws.Range("H103").Value = 1073.2858
$ws.Range("I103").Value = 855
$ws.Range("J103").Value = 1237
$ws.Range("K103").Value = 2565
$ws.Range("L103").Value = 3711
$ws.Range("M103").Value = -1686
$ws.Range("N103").Value = -5469
$ws.Range("H122").Value = 7144096
$ws.Range("J122").Value = 9091800
$ws.Range("L122").Value = 81826200
$ws.Range("N122").Value = -81831100
$ws.Range("H131").Value = 21430.36
$ws.Range("J131").Value = 1591.262
$ws.Range("L131").Value = 4773.786
$ws.Range("N131").Value = -14853.786
$ws.Range("H135").Value = 1985.0968
$ws.Range("I135").Value = 1103.7
$ws.Range("J135").Value = 2404.8096
$ws.Range("K135").Value = 9933.300000000001
$ws.Range("L135").Value = 21643.2864
$ws.Range("M135").Value = -7398.300000000001
$ws.Range("N135").Value = -26713.2864
$ws.Range("H136").Value = 998.5714
$ws.Range("I136").Value = 833
$ws.Range("K136").Value = 2499
$ws.Range("M136").Value = 2601
$ws.Range("H137").Value = 2882.3333
$ws.Range("I137").Value = 2855.8
$ws.Range("J137").Value = 2895.6
$ws.Range("K137").Value = 8567.400000000001
$ws.Range("L137").Value = 8686.799999999999
$ws.Range("M137").Value = -3467.400000000001
$ws.Range("N137").Value = -18886.8
$ws.Range("H138").Value = 2394
$ws.Range("I138").Value = 2394
$ws.Range("K138").Value = 7182
$ws.Range("M138").Value = -2042

# --- Sheet: GSM (18 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2527.2144
$ws.Range("I122").Value = 1837
$ws.Range("J122").Value = 3447.5
$ws.Range("K122").Value = 5511
$ws.Range("L122").Value = 10342.5
$ws.Range("M122").Value = -3061
$ws.Range("N122").Value = -15242.5
$ws.Range("H126").Value = 7010.625
$ws.Range("I126").Value = 6843
$ws.Range("J126").Value = 7513.5
$ws.Range("K126").Value = 20529
$ws.Range("L126").Value = 22540.5
$ws.Range("M126").Value = -18059
$ws.Range("N126").Value = -27480.5
$ws.Range("H132").Value = 3381.4546
$ws.Range("I132").Value = 2577.4443
$ws.Range("K132").Value = 7732.3329
$ws.Range("M132").Value = -5202.3329

# --- Sheet: LTW (22 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2232.125
$ws.Range("I61").Value = 2122.4285
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2122.4285
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1920.4285
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 2232.125
$ws.Range("I113").Value = 2122.4285
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2122.4285
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 47.57150000000001
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 5256.143
$ws.Range("I122").Value = 4938.6
$ws.Range("K122").Value = 14815.8
$ws.Range("M122").Value = -12365.8
$ws.Range("H132").Value = 12348785
$ws.Range("I132").Value = 12348785
$ws.Range("K132").Value = 37046355
$ws.Range("M132").Value = -37043825

# --- Sheet: WVR (15 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2416.1875
$ws.Range("I96").Value = 2722.8333
$ws.Range("K96").Value = 2722.8333
$ws.Range("M96").Value = -1349.8333
$ws.Range("H113").Value = 413.57693
$ws.Range("I113").Value = 356.5625
$ws.Range("J113").Value = 504.8
$ws.Range("K113").Value = 1069.6875
$ws.Range("L113").Value = 1514.4
$ws.Range("M113").Value = 1100.3125
$ws.Range("N113").Value = -5854.4
$ws.Range("H122").Value = 3016.45
$ws.Range("I122").Value = 3016.45
$ws.Range("K122").Value = 9049.349999999999
$ws.Range("M122").Value = -6599.349999999999
